# Updates cryptos list: refresh Price (D) and Volume(1h) (E) columns for rows 2-51.
# Mirrors the data-refresh GitHub Action commit described in the task.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.643.84"
$ws.Range("E2").Value = "  +0.50%  "

# Row 3
$ws.Range("D3").Value = "2.114.61"
$ws.Range("E3").Value = "  +0.15%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.21%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "339.42"
$ws.Range("E5").Value = "  +1.55%  "

# Row 6
$ws.Range("E6").Value = "  +1.02%  "

# Row 7
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4502"
$ws.Range("E8").Value = "  +0.12%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.87"
$ws.Range("E9").Value = "  +0.69%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09036"
$ws.Range("E10").Value = "  -0.31%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.171"
$ws.Range("E11").Value = "  -0.34%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.35"
$ws.Range("E12").Value = "  -0.55%  "

# Row 13
$ws.Range("D13").Value = "2.125.27"
$ws.Range("E13").Value = "  +1.25%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.788"
$ws.Range("E14").Value = "  -0.07%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.059"
$ws.Range("E15").Value = "  +3.04%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "97.83"
$ws.Range("E16").Value = "  +1.07%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001162"
$ws.Range("E17").Value = "  +2.70%  "

# Row 18
$ws.Range("E18").Value = "  +1.08%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06706"
$ws.Range("E19").Value = "  +1.27%  "

# Row 20
$ws.Range("E20").Value = "  -0.38%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.406"
$ws.Range("E22").Value = "  +1.31%  "

# Row 23
$ws.Range("D23").Value = "30.740.13"
$ws.Range("E23").Value = "  +0.66%  "

# Row 24
$ws.Range("E24").Value = "  +3.16%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.378"
$ws.Range("E25").Value = "  +0.90%  "

# Row 26
$ws.Range("D26").Value = "2.373.53"
$ws.Range("E26").Value = "  +1.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.40"
$ws.Range("E27").Value = "  +0.00%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "165.64"
$ws.Range("E28").Value = "  +1.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.557"
$ws.Range("E29").Value = "  -1.48%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "135.62"
$ws.Range("E30").Value = "  +2.02%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.196"
$ws.Range("E31").Value = "  -0.34%  "

# Row 32
$ws.Range("E32").Value = "  +0.01%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.384"
$ws.Range("E33").Value = "  +3.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.627"
$ws.Range("E34").Value = "  -2.49%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.940"
$ws.Range("E35").Value = "  +0.33%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "10.34"
$ws.Range("E36").Value = "  -2.61%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.897"
$ws.Range("E37").Value = "  +5.45%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02655"
$ws.Range("E38").Value = "  +2.79%  "

# Row 39
$ws.Range("E39").Value = "  -0.07%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2325"
$ws.Range("E40").Value = "  +0.57%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.60"
$ws.Range("E41").Value = "  -1.55%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6875"
$ws.Range("E42").Value = "  -0.89%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.260"
$ws.Range("E43").Value = "  +1.00%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.07"
$ws.Range("E44").Value = "  +6.56%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6430"
$ws.Range("E45").Value = "  +0.26%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.315"
$ws.Range("E46").Value = "  -2.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000370"
$ws.Range("E47").Value = "  +11.97%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.709"
$ws.Range("E48").Value = "  +1.14%  "

# Row 49
$ws.Range("E49").Value = "  +0.40%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "82.93"
$ws.Range("E50").Value = "  -0.68%  "

# Row 51
$ws.Range("E51").Value = "  +3.16%  "
